$d = $word.ActiveDocument

$d.Content.Find.Execute("Brimit.Silktide.Umbraco.Plugin", $true, $false, $false, $false, $false, $true, 1, $false, "Silktide.Umbraco.Plugin", 2)
